$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column C (rows 2 through 252) to the uniform value 7310
$ws.Range("C2:C252").Value = 7310
